$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-09-03 01:03:58"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-09-03 01:03:53"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-09-03 01:03:58"

# --- column width changes ---
# Target stored XML width is 17.2159881591797 characters; the engine quantizes
# ColumnWidth to a 1/6-character pixel grid, so use the nearest attainable
# value (16 + 1/3 characters -> stored width 17.1666...).
$newWidth = 16 + (1/3)
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
